$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.221942067146301
$ws.Range("B1").Value = 3.559767484664917
$ws.Range("C1").Value = 2.464468955993652
$ws.Range("D1").Value = 1.287385702133179
$ws.Range("E1").Value = 0.9550935029983521
